$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-22 07:04:14"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-22 07:04:09"
$zhcn.Range("K2").Value = "2016-08-22 07:04:34"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-08-22 07:04:41"
